# superdataset-20 (without cons) test on mae
# Populate the J (test train MAE) / K (test MAE) columns for rows 5-54
# with the captured test-run results for the "without cons" model,
# and scroll/select the view the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$results = @(
    @(5, 55.424267908505684, 154.4179244249726),
    @(6, 55.259116559375421, 144.76450164293539),
    @(7, 55.017023695384189, 154.58291894852141),
    @(8, 55.424318586494998, 145.71533953997809),
    @(9, 56.229327489385021, 143.66248630887179),
    @(10, 55.303647445555399, 146.41512595837901),
    @(11, 55.923851527188063, 142.85158269441399),
    @(12, 55.045093822764009, 147.25348849945229),
    @(13, 56.149997260649229, 141.59311610076671),
    @(14, 55.989398712505128, 141.10587075575029),
    @(15, 54.829087796192297, 150.3983515881709),
    @(16, 54.795221202574993, 153.0842332968237),
    @(17, 55.664584303520073, 135.31889375684551),
    @(18, 55.996336118339947, 142.65653340635271),
    @(19, 55.156718257772901, 149.92729463307779),
    @(20, 55.584656896315558, 142.1763855421687),
    @(21, 54.96366388166004, 147.24496166484121),
    @(22, 54.568407067524987, 155.1026067907996),
    @(23, 54.478741268319411, 154.83286966046001),
    @(24, 55.370594439117923, 148.98651150054761),
    @(25, 55.134302150390347, 155.00518072289151),
    @(26, 55.907507190795783, 131.24786418400879),
    @(27, 55.533096836049857, 146.29264512595839),
    @(28, 55.050604026845633, 147.51521358159911),
    @(29, 54.789508286536091, 153.77193318729459),
    @(30, 54.549936994932203, 148.3796166484118),
    @(31, 55.29320914943159, 148.67573384446879),
    @(32, 55.050753321462807, 146.82174151150051),
    @(33, 55.093511847692092, 140.44404709748079),
    @(34, 54.987641418983713, 152.9521029572837),
    @(35, 54.791792905081493, 151.64943592552021),
    @(36, 54.873843309135729, 154.6143702081051),
    @(37, 55.249631557320903, 152.56055859802851),
    @(38, 54.881038213943278, 142.5526725082147),
    @(39, 56.051006711409393, 137.52119934282581),
    @(40, 54.176366251198459, 157.27287513691121),
    @(41, 53.988820709491847, 163.01746440306681),
    @(42, 56.263284481577863, 142.6547864184009),
    @(43, 56.314307629091907, 138.63051478641839),
    @(44, 56.142250376660719, 139.37358159912381),
    @(45, 55.18708807012738, 150.22291894852131),
    @(46, 55.127571565538958, 150.3330777656079),
    @(47, 54.198056430625932, 149.42481380065709),
    @(48, 55.338035885495131, 145.1532584884994),
    @(49, 55.952420216408697, 145.642004381161),
    @(50, 55.455764963703608, 145.42623767798469),
    @(51, 55.910698534447327, 139.0867743702081),
    @(52, 54.496575811532658, 146.08871851040519),
    @(53, 54.544969182303767, 148.86469331872951),
    @(54, 55.698081084782899, 142.90341182913471),
)

foreach ($item in $results) {
    $row = $item[0]
    $ws.Cells.Item($row, 10).Value = $item[1]   # column J - train (MAE)
    $ws.Cells.Item($row, 11).Value = $item[2]   # column K - test (MAE)
}

# Restore the view: scrolled so row 19 is at the top, with L41 selected
$win = $excel.ActiveWindow
$ws.Range("L41").Select()
$win.ScrollRow = 19
$win.ScrollColumn = 1
